$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B = Coin, C = Link, D = Price, E = Volume(1h)
# Each row entry only lists the columns that actually change.
$updates = @{
    2  = @{ D = "42.466.04"; E = "  +1.62%  " }
    3  = @{ D = "2.281.37";  E = "  +0.46%  " }
    4  = @{ E = "  +0.03%  " }
    5  = @{ D = "307.34";    E = "  +1.26%  " }
    6  = @{ D = "97.65";     E = "  +5.47%  " }
    7  = @{ D = "0.530";     E = "  +0.08%  " }
    8  = @{ E = "  +0.00%  " }
    9  = @{ D = "0.494";     E = "  +1.98%  " }
    10 = @{ D = "36.00";     E = "  +10.62%  " }
    11 = @{ D = "0.0798";    E = "  +0.13%  " }
    12 = @{ E = "  -1.13%  " }
    13 = @{ D = "6.70";      E = "  +0.13%  " }
    14 = @{ D = "2.638.59";  E = "  +0.65%  " }
    15 = @{ D = "14.43";     E = "  +1.01%  " }
    16 = @{ D = "2.283.76";  E = "  -0.29%  " }
    17 = @{ D = "0.797";     E = "  +2.51%  " }
    18 = @{ D = "42.373.82"; E = "  +1.62%  " }
    19 = @{ D = "12.53";     E = "  +0.64%  " }
    20 = @{ D = "0.0₃0912";  E = "  +0.76%  " }
    21 = @{ D = "5.97";      E = "  +0.46%  " }
    22 = @{ D = "67.80";     E = "  +1.00%  " }
    23 = @{ D = "241.22";    E = "  +0.51%  " }
    24 = @{ D = "2.60";      E = "  +0.65%  " }
    25 = @{ E = "  +1.52%  " }
    26 = @{ E = "  +0.01%  " }
    27 = @{ D = "23.90";     E = "  -0.36%  " }
    28 = @{ D = "37.61";     E = "  +5.70%  " }
    29 = @{ D = "9.53";      E = "  -0.18%  " }
    30 = @{ E = "  +1.99%  " }
    31 = @{ D = "159.26";    E = "  -0.95%  " }
    32 = @{ D = "5.26";      E = "  +0.26%  " }
    33 = @{ D = "1.00";      E = "  +0.05%  " }
    34 = @{ D = "3.14";      E = "  +4.43%  " }
    35 = @{ D = "0.0742";    E = "  -0.38%  " }
    36 = @{ D = "17.07" }
    37 = @{ B = "WEMIXToken"; C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D = "2.38";  E = "  +0.26%  " }
    38 = @{ B = "Kaspa";      C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas";        D = "0.106"; E = "  +0.77%  " }
    39 = @{ D = "1.84";      E = "  +2.59%  " }
    40 = @{ E = "  -1.18%  " }
    41 = @{ D = "4.10";      E = "  +4.65%  " }
    42 = @{ E = "  +14.62%  " }
    43 = @{ D = "2.000.27";  E = "  -0.11%  " }
    44 = @{ D = "0.0285";    E = "  +1.46%  " }
    45 = @{ D = "18.82";     E = "  -1.83%  " }
    46 = @{ D = "2.96";      E = "  +1.77%  " }
    47 = @{ D = "10.00";     E = "  -3.34%  " }
    48 = @{ D = "52.95";     E = "  +0.58%  " }
    49 = @{ E = "  +0.84%  " }
    50 = @{ D = "72.20";     E = "  +0.24%  " }
    51 = @{ D = "92.19";     E = "  +1.12%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        if ($col -eq "D") {
            # Price column holds text that often looks numeric (e.g. "307.34").
            # Force it to be written back as text, matching the source data,
            # then drop back to the Normal style so no stray number format
            # is left attached to the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $cols[$col]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $cols[$col]
        }
    }
}
